# Update Release-Notes.xlsx - Folder inventory updated on Thu Jun 12 16:20:33 UTC 2025
$wb = $excel.ActiveWorkbook

# --- Sheet: Folder Inventory ---
$wsInventory = $wb.Worksheets.Item("Folder Inventory")
$wsInventory.Range("C2").Value = "2025-06-12 21:50:14 +0530"

# --- Sheet: Metadata ---
$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B3").Value = "2025-06-12 16:20:33 UTC"
$wsMetadata.Range("B5").NumberFormat = "@"
$wsMetadata.Range("B5").Value = "7"

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = "2025-06-12 21:50:14 +0530"
